$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure these cells stay text (some "new" values look numeric, e.g. "0.999", "7.69")
# so Excel must not auto-convert them to numbers - preserves exact formatting/leading/trailing chars.
$ws.Range("B2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '61.187.55'
$ws.Range("E2").Value = '  -1.22%  '
$ws.Range("D3").Value = '3.389.00'
$ws.Range("E3").Value = '  -0.26%  '
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.17%  '
$ws.Range("D5").Value = '575.82'
$ws.Range("E5").Value = '  -0.84%  '
$ws.Range("D6").Value = '137.77'
$ws.Range("E6").Value = '  -0.48%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").Value = '3.385.79'
$ws.Range("E8").Value = '  -0.33%  '
$ws.Range("D9").Value = '0.470'
$ws.Range("E9").Value = '  -1.05%  '
$ws.Range("D10").Value = '7.69'
$ws.Range("E10").Value = '  +2.27%  '
$ws.Range("D11").Value = '0.123'
$ws.Range("E11").Value = '  -3.30%  '
$ws.Range("D12").Value = '0.383'
$ws.Range("E12").Value = '  -2.37%  '
$ws.Range("D13").Value = '3.963.65'
$ws.Range("E13").Value = '  -0.51%  '
$ws.Range("E14").Value = '  +0.93%  '
$ws.Range("D15").Value = '0.0000173'
$ws.Range("E15").Value = '  -3.23%  '
$ws.Range("D16").Value = '3.387.14'
$ws.Range("E16").Value = '  -0.23%  '
$ws.Range("D17").Value = '25.65'
$ws.Range("E17").Value = '  +0.79%  '
$ws.Range("D18").Value = '61.265.34'
$ws.Range("E18").Value = '  -1.23%  '
$ws.Range("D19").Value = '13.89'
$ws.Range("E19").Value = '  -2.08%  '
$ws.Range("B20").Value = 'Uniswap'
$ws.Range("C20").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D20").Value = '9.40'
$ws.Range("E20").Value = '  -1.11%  '
$ws.Range("B21").Value = 'Polkadot'
$ws.Range("C21").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D21").Value = '5.76'
$ws.Range("E21").Value = '  -1.20%  '
$ws.Range("D22").Value = '379.21'
$ws.Range("E22").Value = '  -3.85%  '
$ws.Range("B23").Value = 'Polygon'
$ws.Range("C23").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D23").Value = '0.552'
$ws.Range("E23").Value = '  -2.49%  '
$ws.Range("B24").Value = 'WrappedeETH'
$ws.Range("C24").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D24").Value = '3.522.39'
$ws.Range("E24").Value = '  -0.70%  '
$ws.Range("E25").Value = '  +0.39%  '
$ws.Range("D26").Value = '0.0000127'
$ws.Range("E26").Value = '  -2.41%  '
$ws.Range("D27").Value = '71.16'
$ws.Range("E27").Value = '  -0.62%  '
$ws.Range("D28").Value = '0.182'
$ws.Range("E28").Value = '  +12.74%  '
$ws.Range("D29").Value = '1.66'
$ws.Range("E29").Value = '  +0.05%  '
$ws.Range("D30").Value = '0.999'
$ws.Range("E30").Value = '  -0.04%  '
$ws.Range("D31").Value = '7.42'
$ws.Range("E31").Value = '  -3.20%  '
$ws.Range("D32").Value = '8.11'
$ws.Range("E32").Value = '  -1.56%  '
$ws.Range("D33").Value = '2.15'
$ws.Range("E33").Value = '  -1.49%  '
$ws.Range("E34").Value = '  -0.03%  '
$ws.Range("D35").Value = '23.51'
$ws.Range("E35").Value = '  -0.01%  '
$ws.Range("E36").Value = '  -4.45%  '
$ws.Range("D37").Value = '1.55'
$ws.Range("E37").Value = '  -3.03%  '
$ws.Range("D38").Value = '6.85'
$ws.Range("E38").Value = '  -1.03%  '
$ws.Range("D39").Value = '164.44'
$ws.Range("E39").Value = '  -0.31%  '
$ws.Range("D40").Value = '0.0759'
$ws.Range("E40").Value = '  -3.87%  '
$ws.Range("B41").Value = 'EnergySwap'
$ws.Range("C41").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D41").Value = '25.70'
$ws.Range("E41").Value = '  +2.95%  '
$ws.Range("B42").Value = 'FirstDigitalUSD'
$ws.Range("C42").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D42").Value = '0.998'
$ws.Range("E42").Value = '  -0.22%  '
$ws.Range("D43").Value = '0.775'
$ws.Range("E43").Value = '  -1.70%  '
$ws.Range("B44").Value = 'Stacks'
$ws.Range("C44").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D44").Value = '1.71'
$ws.Range("E44").Value = '  -4.06%  '
$ws.Range("B45").Value = 'OKB'
$ws.Range("C45").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D45").Value = '41.71'
$ws.Range("E45").Value = '  +0.88%  '
$ws.Range("D46").Value = '4.36'
$ws.Range("E46").Value = '  -1.80%  '
$ws.Range("E47").Value = '  -5.09%  '
$ws.Range("D48").Value = '2.518.09'
$ws.Range("E48").Value = '  +7.53%  '
$ws.Range("D49").Value = '6.80'
$ws.Range("E49").Value = '  -1.64%  '
$ws.Range("D50").Value = '22.94'
$ws.Range("E50").Value = '  -1.05%  '
$ws.Range("D51").Value = '2.47'
$ws.Range("E51").Value = '  +5.33%  '
